# Workbook: Swag_lab_data.xlsx
# Commit: "impledted new pages and test"
#
# Adds two new worksheets - "verifycheckoutOverview" and
# "verifyFinishCheckout" - right after "verifyCancelBtn" and before
# "invalidLoginTest", each holding the same login/checkout-form test data
# (username, password, firstname, lastname, zipcode) used elsewhere in the
# workbook, plus a new zipcode value "R42201".

$wb = $excel.ActiveWorkbook

$anchor = $wb.Worksheets.Item("verifyCancelBtn")

# --- New sheet 1: verifycheckoutOverview -----------------------------
$overview = $wb.Worksheets.Add($null, $anchor)
$overview.Name = "verifycheckoutOverview"

$overview.Cells.Item(1, 1).Value = "username"
$overview.Cells.Item(1, 2).Value = "password"
$overview.Cells.Item(1, 3).Value = "firstname"
$overview.Cells.Item(1, 4).Value = "lastname"
$overview.Cells.Item(1, 5).Value = "zipcode"

$overview.Cells.Item(2, 1).Value = "standard_user"
$overview.Cells.Item(2, 2).Value = "secret_sauce"
$overview.Cells.Item(2, 3).Value = "Rahul"
$overview.Cells.Item(2, 4).Value = "Kashyap"
$overview.Cells.Item(2, 5).Value = "R42201"

$overview.Range("G9").Select() | Out-Null

# --- New sheet 2: verifyFinishCheckout --------------------------------
$finish = $wb.Worksheets.Add($null, $overview)
$finish.Name = "verifyFinishCheckout"

$finish.Cells.Item(1, 1).Value = "username"
$finish.Cells.Item(1, 2).Value = "password"
$finish.Cells.Item(1, 3).Value = "firstname"
$finish.Cells.Item(1, 4).Value = "lastname"
$finish.Cells.Item(1, 5).Value = "zipcode"

$finish.Cells.Item(2, 1).Value = "standard_user"
$finish.Cells.Item(2, 2).Value = "secret_sauce"
$finish.Cells.Item(2, 3).Value = "Rahul"
$finish.Cells.Item(2, 4).Value = "Kashyap"
$finish.Cells.Item(2, 5).Value = "R42201"

$finish.Range("F5").Select() | Out-Null

# verifyFinishCheckout ends up the active/selected tab, matching the
# workbook's new activeTab.
$finish.Activate() | Out-Null
